$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Before: row5 = 2017-09-26 / 5000 / test2 / 60   (top-border style, right under header)
#         row6 = (blank date) / 500000 / - / (J blank) (no-border style)
#
# After:  row5 = 2017-09-24 / 150 / -                      (NEW, top-border style)
#         row6 = 2017-09-26 / 5000 / test2 / 60            (old row5 data, now no-border style)
#         row7 = (blank date) / 500000 / -                 (old row6 data, no-border style, IN-KIND merge widened)
#         row8 = 2017-09-27 / - / test3 / 58               (NEW, no-border style)
# -----------------------------------------------------------------

# 1. Insert a new blank row at position 5; this pushes the existing
#    row5 -> row6 and row6 -> row7 (values, formats and merges all shift).
$ws.Rows.Item(5).Insert()

# 2. Build the new row 5 (top-border style, same pattern the old row5 had -
#    which is now sitting, unchanged, on row 6 after the shift above).
$ws.Range("E5:F5").Merge()

$ws.Range("A6:F6").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)

$ws.Range("G6").Copy()
$ws.Range("G5:J5").PasteSpecial(-4122)
$ws.Range("G5:J5").HorizontalAlignment = -4108

$ws.Range("A5:D5").Merge()
$ws.Range("G5:J5").Merge()

$ws.Range("A5").Value = 43002
$ws.Range("E5").Value = 150
$ws.Range("G5").Value = "-"

# 3. Row 6 (shifted old row5 data) must switch from the "top-border" style
#    it inherited to the "no-border" style, since it is no longer the row
#    directly under the header. Pull that style from row 7 (shifted old
#    row6 data), which already has the correct no-border look.
$ws.Range("A7:F7").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)

$ws.Range("G7").Copy()
$ws.Range("G6:I6").PasteSpecial(-4122)
$ws.Range("G6:I6").HorizontalAlignment = -4108

# 4. Row 7 (shifted old row6 data): clear the date cells (A7:D7 had no
#    value in the target) and widen the IN-KIND merge from G:I to G:J.
$ws.Range("A7:D7").ClearContents()
$ws.Range("G7:I7").UnMerge()
$ws.Range("G7:J7").Merge()

# 5. Build brand-new row 8 (no-border style, same pattern as row 6/7).
$ws.Range("A8:D8").Merge()
$ws.Range("E8:F8").Merge()
$ws.Range("G8:I8").Merge()

$ws.Range("A7:F7").Copy()
$ws.Range("A8:F8").PasteSpecial(-4122)

$ws.Range("G6:I6").Copy()
$ws.Range("G8:I8").PasteSpecial(-4122)

$ws.Range("A8").Value = 43005
$ws.Range("E8").Value = "-"
$ws.Range("G8").Value = "test3"
$ws.Range("J8").Value = 58
